# DDAf_2023_Tableau_annexe_sources.xlsx — refresh "last updated" source dates
# and correct the World Urbanization Prospects reference (year 2022 -> 2018,
# plus its own "2018" reference year in the "Derniere mise a jour" column).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sources")

# Comite d'aide au developpement de l'OCDE (CAD) - row 4
$ws.Range("D4").Value = "Récupéré le 27/10/2023"

# FMI Balance of Payments and IIP Statistics - row 13
$ws.Range("D13").Value = "Mis à jour le 25/10/2023"

# Global Knowledge Partnership on Migration and Development (KNOMAD) - row 16
$ws.Range("D16").Value = "Juin 2023"

# Indicateurs de developpement de la Banque mondiale (WDI) - row 19
$ws.Range("D19").Value = "Mise à jour le 10/10/2023"

# Perspectives de l'economie mondiale du FMI (WEO) - row 26
$ws.Range("D26").Value = "Octobre 2023"

# UNCTADSTAT, Centre de donnees sur les IDE - row 28
$ws.Range("D28").Value = "Mis à jour le 22/09/2022"

# World Urbanization Prospects - row 33: title year 2022 -> 2018,
# and the update-date cell now holds the text "2018" instead of numeric 2021
$ws.Range("B33").Value = "World Urbanization Prospects 2018, ONU"
$ws.Range("D33").Value = "2018"
